$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values (rows 2-10, columns E:T)
# reflects rerun of NATMI lrc2p computation with new TPM input

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7170026666666667
$ws.Range("H2").Value = 2.151008
$ws.Range("I2").Value = 0.02953485643833859
$ws.Range("J2").Value = 0.02953485643833859
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 15.75563966666667
$ws.Range("N2").Value = 47.266919
$ws.Range("O2").Value = 0.3220556913988901
$ws.Range("P2").Value = 0.32205569139889
$ws.Range("Q2").Value = 11.29683565603911
$ws.Range("R2").Value = 101.671520904352
$ws.Range("S2").Value = 0.009511868610616094
$ws.Range("T2").Value = 0.009511868610616094

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7170026666666667
$ws.Range("H3").Value = 2.151008
$ws.Range("I3").Value = 0.02953485643833859
$ws.Range("J3").Value = 0.02953485643833859
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 26.95384733333333
$ws.Range("N3").Value = 80.861542
$ws.Range("O3").Value = 0.5509544596378365
$ws.Range("P3").Value = 0.5509544596378364
$ws.Range("Q3").Value = 19.32598041492622
$ws.Range("R3").Value = 173.933823734336
$ws.Range("S3").Value = 0.01627236086946592
$ws.Range("T3").Value = 0.01627236086946591

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7170026666666667
$ws.Range("H4").Value = 2.151008
$ws.Range("I4").Value = 0.02953485643833859
$ws.Range("J4").Value = 0.02953485643833859
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 6.212609666666666
$ws.Range("N4").Value = 18.637829
$ws.Range("O4").Value = 0.1269898489632735
$ws.Range("P4").Value = 0.1269898489632735
$ws.Range("Q4").Value = 4.454457697959111
$ws.Range("R4").Value = 40.090119281632
$ws.Range("S4").Value = 0.003750626958256584
$ws.Range("T4").Value = 0.003750626958256584

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 12.06704733333333
$ws.Range("H5").Value = 36.201142
$ws.Range("I5").Value = 0.4970672037825566
$ws.Range("J5").Value = 0.4970672037825566
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 15.75563966666667
$ws.Range("N5").Value = 47.266919
$ws.Range("O5").Value = 0.3220556913988901
$ws.Range("P5").Value = 0.32205569139889
$ws.Range("Q5").Value = 190.1240496246109
$ws.Range("R5").Value = 1711.116446621498
$ws.Range("S5").Value = 0.1600833219859042
$ws.Range("T5").Value = 0.1600833219859042

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 12.06704733333333
$ws.Range("H6").Value = 36.201142
$ws.Range("I6").Value = 0.4970672037825566
$ws.Range("J6").Value = 0.4970672037825566
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 26.95384733333333
$ws.Range("N6").Value = 80.861542
$ws.Range("O6").Value = 0.5509544596378365
$ws.Range("P6").Value = 0.5509544596378364
$ws.Range("Q6").Value = 325.2533515867738
$ws.Range("R6").Value = 2927.280164280964
$ws.Range("S6").Value = 0.2738613926637088
$ws.Range("T6").Value = 0.2738613926637088

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 12.06704733333333
$ws.Range("H7").Value = 36.201142
$ws.Range("I7").Value = 0.4970672037825566
$ws.Range("J7").Value = 0.4970672037825566
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 6.212609666666666
$ws.Range("N7").Value = 18.637829
$ws.Range("O7").Value = 0.1269898489632735
$ws.Range("P7").Value = 0.1269898489632735
$ws.Range("Q7").Value = 74.96785491119088
$ws.Range("R7").Value = 674.710694200718
$ws.Range("S7").Value = 0.06312248913294356
$ws.Range("T7").Value = 0.06312248913294355

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 11.49244066666667
$ws.Range("H8").Value = 34.477322
$ws.Range("I8").Value = 0.4733979397791048
$ws.Range("J8").Value = 0.4733979397791048
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 15.75563966666667
$ws.Range("N8").Value = 47.266919
$ws.Range("O8").Value = 0.3220556913988901
$ws.Range("P8").Value = 0.32205569139889
$ws.Range("Q8").Value = 181.0707540345465
$ws.Range("R8").Value = 1629.636786310918
$ws.Range("S8").Value = 0.1524605008023697
$ws.Range("T8").Value = 0.1524605008023697

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 11.49244066666667
$ws.Range("H9").Value = 34.477322
$ws.Range("I9").Value = 0.4733979397791048
$ws.Range("J9").Value = 0.4733979397791048
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 26.95384733333333
$ws.Range("N9").Value = 80.861542
$ws.Range("O9").Value = 0.5509544596378365
$ws.Range("P9").Value = 0.5509544596378364
$ws.Range("Q9").Value = 309.7654912167249
$ws.Range("R9").Value = 2787.889420950524
$ws.Range("S9").Value = 0.2608207061046618
$ws.Range("T9").Value = 0.2608207061046617

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 11.49244066666667
$ws.Range("H10").Value = 34.477322
$ws.Range("I10").Value = 0.4733979397791048
$ws.Range("J10").Value = 0.4733979397791048
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 6.212609666666666
$ws.Range("N10").Value = 18.637829
$ws.Range("O10").Value = 0.1269898489632735
$ws.Range("P10").Value = 0.1269898489632735
$ws.Range("Q10").Value = 71.39804797932645
$ws.Range("R10").Value = 642.582431813938
$ws.Range("S10").Value = 0.06011673287207338
$ws.Range("T10").Value = 0.06011673287207337
